# "updated reject and withdraw files"
# Sheet1 column A holds a list of IDs; refresh it with the new export:
# replace the existing 10 rows and append 3 new ones (A1:A13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    710006,
    710007,
    710014,
    710019,
    710034,
    710077,
    710083,
    710084,
    710003,
    710000,
    710122,
    709877,
    711545
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Leave the cursor where the author left it when they saved the file.
$ws.Range("B4").Select()
